$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC01_02_03")

# Row 2: remove Dealer Suspense from G2; add 1232/1212 to I2/J2 (as text)
$ws.Range("G2").Value = $null
$ws.Range("I2").Value = "'1232"
$ws.Range("J2").Value = "'1212"

# Row 3: add 'random' to A3; clear C3, E3, F3
$ws.Range("A3").Value = "random"
$ws.Range("C3").Value = $null
$ws.Range("E3").Value = $null
$ws.Range("F3").Value = $null

# Row 4: clear I4, J4 (1233/1213 removed)
$ws.Range("I4").Value = $null
$ws.Range("J4").Value = $null

# Update selection to G2 as per sheetView change
$ws.Range("G2").Select()
